$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.122.35'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.00%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.892.51'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7397'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -1.75%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '242.69'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.15%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9996'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3174'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +1.69%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07213'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +1.24%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.92'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -1.54%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08333'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -1.80%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.017.19'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +7.55%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.7604'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.440'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.54%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '93.24'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.15%  '
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.53%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '30.243.03'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.29%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '251.05'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.43%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007896'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.19%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.210.49'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +2.94%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.003'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.38%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.946'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.86%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.000'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.07%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1578'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -1.14%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.319'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.50%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '164.60'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +1.44%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.77'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.19%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.067'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +1.90%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.42%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.595'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +1.87%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.540'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.04%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.205'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +1.84%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05369'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -1.00%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +1.38%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7760'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +3.50%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.45%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.724'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.63%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01962'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.62%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.766'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.23%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4572'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +2.53%  '
$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.110'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.01%  '
$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.099.25'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.85%  '
$ws.Range('B44').NumberFormat = '@'
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').NumberFormat = '@'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '72.74'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.14%  '
$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value = 'TrustWalletToken'
$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8822'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +2.86%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '104.54'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +2.02%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.000'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.05%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.867'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.46%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.598'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -1.64%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.123.96'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +4.54%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '9.624'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -1.43%  '
